$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Customer Quote")

# Fix surplus number: change K16, K18, K21, K22, K26 from 1.0565 to 1
$ws.Range("K16").Value = 1
$ws.Range("K18").Value = 1
$ws.Range("K21").Value = 1
$ws.Range("K22").Value = 1
$ws.Range("K26").Value = 1

# Update active cell selection to H6 (support for longer quotes view)
$ws.Range("H6").Select()
